$d = $word.ActiveDocument

# Replace the hard-coded legal representative name with the ${represent} placeholder.
$d.Content.Find.Execute(
    "NGUYỄN HUỲNH THU TRÚC", $false, $false, $false, $false, $false,
    $true, 1, $false, "`${represent}", 2)

# Replace the hard-coded "Giám Đốc" title with the ${position} placeholder and
# append ", làm đại diện" after it.
$d.Content.Find.Execute(
    "Chức vụ: Giám Đốc", $false, $false, $false, $false, $false,
    $true, 1, $false, "Chức vụ: `${position}, làm đại diện", 2)
